$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# colors (COLORREF / BGR order, as consumed by this host)
$orange = 49407   # FFC000
$darkred = 192     # C00000

# --- Table titles ---
$ws.Range("B3").Value = "dojos"
$ws.Range("F3").Value = "ubicaciónes"
$ws.Range("J3").Value = "ninjas"

# --- dojos table headers (B4:D4) ---
$ws.Range("B4").Value = "id"
$ws.Range("B4").Interior.Color = $orange
$ws.Range("C4").Value = "calle"
$ws.Range("D4").Value = "id_ubicacion"
$ws.Range("D4").Interior.Color = $darkred

# --- ubicaciónes table headers (F4:H4) ---
$ws.Range("F4").Value = "id"
$ws.Range("F4").Interior.Color = $orange
$ws.Range("G4").Value = "nombre"
$ws.Range("H4").Value = "numero"

# --- ninjas table headers (J4:L4), M4 column removed ---
$ws.Range("J4").Value = "id"
$ws.Range("J4").Interior.Color = $orange
$ws.Range("K4").Value = "calle"
$ws.Range("L4").Value = "id_dojo"
$ws.Range("L4").Interior.Color = $darkred
$ws.Range("M4:M5").ClearContents()

# --- data rows ---
$ws.Range("G5").Value = "rosas"

$ws.Range("F6").Value = 2
$ws.Range("G6").Value = "ciruelos"
$ws.Range("H6").Value = 33
$ws.Range("K6").Value = "gory"
$ws.Range("L6").Value = 2

# --- habilidades table title + habilidades_ninjas table title ---
$ws.Range("F11").Value = "habilidades"
$ws.Range("J11").Value = "habilidades_ninjas"

# --- habilidades table headers (F12:G12) ---
$ws.Range("F12").Value = "id"
$ws.Range("F12").Interior.Color = $orange

# --- habilidades_ninjas table headers (J12:L12) ---
$ws.Range("J12").Value = "id"
$ws.Range("J12").Interior.Color = $orange
$ws.Range("K12").Value = "id_ninja"
$ws.Range("K12").Interior.Color = $darkred
$ws.Range("L12").Value = "id_habilidad"
$ws.Range("L12").Interior.Color = $darkred

# --- habilidades_ninjas data rows ---
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 3

$ws.Range("F14").Value = 2
$ws.Range("G14").Value = "defender"
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 1

$ws.Range("F15").Value = 3
$ws.Range("G15").Value = "atacar"
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 2

$ws.Range("J16").Value = 4
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 3

# --- column width for J (bestFit-style autosize) ---
$ws.Columns("J").ColumnWidth = 15.5

# --- final selection ---
$ws.Range("J16").Select() | Out-Null
